$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Pin Name (column D) corrections ---
$dUpdates = @{
    56 = "CTRLR Data in <1>"
    57 = "CTRLR Data in <0>"
    58 = "CLK PULSE"
    59 = "NC"
    61 = "CLK"
    62 = "STATE BIT"
    63 = "AUDIO DAC <0>"
    64 = "AUDIO DAC <1>"
    65 = "AUDIO DAC <2>"
    66 = "AUDIO DAC <3>"
    67 = "AUDIO_ROM_SPI_MISO"
    68 = "AUDIO_ROM_SPI_MOSI"
    69 = "AUDIO_ROM_SPI_CLK"
    70 = "AUDIO_ROM_SPI_CE"
    71 = "DOUT <15>"
    72 = "DOUT <14>"
    73 = "DOUT <13>"
    74 = "DOUT <12>"
    75 = "DOUT <11>"
    76 = "DOUT <10>"
    78 = "DOUT <9>"
    79 = "DOUT <8>"
    80 = "DOUT <7>"
    81 = "DOUT <6>"
    82 = "DOUT <5>"
    83 = "DOUT <4>"
    84 = "DOUT <3>"
    85 = "DOUT <2>"
    86 = "DOUT <1>"
    87 = "DOUT <0>"
}
foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# --- Populate Description (column F) ---
$fUpdates = @{
    4 = "SRAM Write Enable"
    5 = "SRAM Output Enabale"
    6 = "SRAM Chip Enable"
    7 = "Address Bit-15 connected to both SRAM and ROM"
    8 = "Address Bit-14 connected to both SRAM and ROM"
    9 = "Address Bit-13 connected to both SRAM and ROM"
    10 = "Ground"
    11 = "Address Bit-12 connected to both SRAM and ROM"
    12 = "Address Bit-11 connected to both SRAM and ROM"
    13 = "Address Bit-10 connected to both SRAM and ROM"
    14 = "Address Bit-9 connected to both SRAM and ROM"
    15 = "Address Bit-8 connected to both SRAM and ROM"
    16 = "Address Bit-7 connected to both SRAM and ROM"
    17 = "Address Bit-6 connected to both SRAM and ROM"
    18 = "Address Bit-5 connected to both SRAM and ROM"
    19 = "Digital Power Supply"
    20 = "Address Bit-4 connected to both SRAM and ROM"
    21 = "Address Bit-3 connected to both SRAM and ROM"
    22 = "Address Bit-2 connected to both SRAM and ROM"
    23 = "Address Bit-1 connected to both SRAM and ROM"
    24 = "Address Bit-0 connected to both SRAM and ROM"
    25 = "Data in Bit-0 connect to both SRAM and ROM"
    26 = "Data in Bit-1 connect to both SRAM and ROM"
    27 = "Data in Bit-2 connect to both SRAM and ROM"
    28 = "Data in Bit-3 connect to both SRAM and ROM"
    29 = "Data in Bit-4 connect to both SRAM and ROM"
    30 = "Data in Bit-5 connect to both SRAM and ROM"
    31 = "Data in Bit-6 connect to both SRAM and ROM"
    32 = "Data in Bit-7 connect to both SRAM and ROM"
    33 = "Data in Bit-8 connect to both SRAM and ROM"
    34 = "Data in Bit-9 connect to both SRAM and ROM"
    35 = "Data in Bit-10 connect to both SRAM and ROM"
    36 = "Digital Power Supply"
    37 = "Data in Bit-11 connect to both SRAM and ROM"
    38 = "Data in Bit-12 connect to both SRAM and ROM"
    39 = "Data in Bit-13 connect to both SRAM and ROM"
    40 = "Data in Bit-14 connect to both SRAM and ROM"
    41 = "Data in Bit-15 connect to both SRAM and ROM"
    42 = "ROM Chip Enable"
    43 = "ROM Output Enable"
    44 = "Least Significant Bit for VGA Color"
    45 = "Most Significant Bit for VGA Color"
    46 = "Horizontal Sync required for VGA"
    47 = "Vertical Sync required for VGA"
    48 = "Program Status Register Bit-0 used for Debugging"
    49 = "Program Status Register Bit-1 used for Debugging"
    50 = "Program Status Register Bit-2 used for Debugging"
    51 = "Ground"
    52 = "Program Status Register Bit-3 used for Debugging"
    53 = "Program Status Register Bit-4 used for Debugging"
    54 = "Controler 2 Latch signal"
    55 = "Controler 1 Latch signal"
    56 = "Controler 2 Data Signal"
    57 = "Controler 1 Data Signal"
    58 = "Clock Output to Both Controllers"
    59 = "Not Connected"
    60 = "Digital Power Supply"
    61 = "Clock Input"
    62 = "Program State Bit used for Debugging"
    63 = "Audio DAC Output Bit-0"
    64 = "Audio DAC Output Bit-1"
    65 = "Audio DAC Output Bit-2"
    66 = "Audio DAC Output Bit-3"
    67 = "Audio ROM SPI Master in Slave out"
    68 = "Audio ROM SPI Master out Slave in"
    69 = "Audio rom SPI Clock"
    70 = "Audio Rom SPI Chip Enable"
    71 = "Data out Bit-15 Connected to SRAM Only"
    72 = "Data out Bit-14 Connected to SRAM Only"
    73 = "Data out Bit-13 Connected to SRAM Only"
    74 = "Data out Bit-12 Connected to SRAM Only"
    75 = "Data out Bit-11 Connected to SRAM Only"
    76 = "Data out Bit-10 Connected to SRAM Only"
    77 = "Ground"
    78 = "Data out Bit-9 Connected to SRAM Only"
    79 = "Data out Bit-8 Connected to SRAM Only"
    80 = "Data out Bit-7 Connected to SRAM Only"
    81 = "Data out Bit-6 Connected to SRAM Only"
    82 = "Data out Bit-5 Connected to SRAM Only"
    83 = "Data out Bit-4 Connected to SRAM Only"
    84 = "Data out Bit-3 Connected to SRAM Only"
    85 = "Data out Bit-2 Connected to SRAM Only"
    86 = "Data out Bit-1 Connected to SRAM Only"
    87 = "Data out Bit-0 Connected to SRAM Only"
}
foreach ($row in $fUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
}

# --- Column width for Description column ---
$ws.Columns.Item(6).ColumnWidth = 42

# --- Reset view / selection ---
$ws.Range("D3:E3").Select()

Write-Host "Pin table updated: $($dUpdates.Count) names corrected, $($fUpdates.Count) descriptions added."
